# Updated cryptos list values (prices + 1h volume %) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.738.51"
$ws.Range("E2").Value = "  +6.44%  "
$ws.Range("D3").Value = "3.474.05"
$ws.Range("E3").Value = "  +4.88%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.42"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  +19.20%  "
$ws.Range("D7").Value = "3.474.73"
$ws.Range("E7").Value = "  +5.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.693"
$ws.Range("E10").Value = "  +8.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.128"
$ws.Range("E11").Value = "  +29.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.13"
$ws.Range("E12").Value = "  +8.53%  "
$ws.Range("D14").Value = "4.014.75"
$ws.Range("E14").Value = "  +4.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.79"
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.20"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "3.473.53"
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").Value = "62.646.26"
$ws.Range("E18").Value = "  +6.80%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.89"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  +25.12%  "
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.34"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.54"
$ws.Range("E24").Value = "  +10.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "315.73"
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.20"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.74"
$ws.Range("E27").Value = "  +8.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.13"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.73"
$ws.Range("E29").Value = "  +4.95%  "
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.38"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "43.86"
$ws.Range("E33").Value = "  +9.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  +23.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.84"
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0497"
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.58"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -9.63%  "
$ws.Range("E42").Value = "  +7.36%  "
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "137.31"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.37"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.02"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.288"
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.31"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "2.234.37"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "3.815.05"
$ws.Range("E51").Value = "  +5.05%  "
